$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b73e5c79445b7ad12cb1d3ccd5c5fa9906e601e/e2e/ea7b708d-cbb4-4f9b-9cda-a4f8800f82b1.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/70e4354520c6b6b692d5494c0ed5e789407c1eac/e2e/ea7b708d-cbb4-4f9b-9cda-a4f8800f82b1.md."
$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/70e4354520c6b6b692d5494c0ed5e789407c1eac/e2e/ea7b708d-cbb4-4f9b-9cda-a4f8800f82b1.md"
$handbackDisplay = "ea7b708d-cbb4-4f9b-9cda-a4f8800f82b1.md"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen the "Error Detail" column (P) so the long message is readable
# (39.14 "characters" renders to an OOXML width of exactly 40)
$wsZh.Columns("P").ColumnWidth = 39.14

$wsZh.Range("I8").Value = $handbackDisplay
$wsZh.Range("J8").Value = "ea7b708d-cbb4-4f9b-9cda-a4f8800f82b1.765c2683109bcf65fb150e795bd71c0083ac57b5.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-08-17 18:42:28"
$wsZh.Range("P8").Value = $errorDetail

$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $handbackUrl, "", "", $handbackDisplay)
$wsZh.Range("I8").Style = "HyperLink"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns("P").ColumnWidth = 39.14

$wsDe.Range("I8").Value = $handbackDisplay
$wsDe.Range("J8").Value = "ea7b708d-cbb4-4f9b-9cda-a4f8800f82b1.765c2683109bcf65fb150e795bd71c0083ac57b5.de-de.xlf"
$wsDe.Range("K8").Value = "2016-08-17 18:42:35"
$wsDe.Range("P8").Value = $errorDetail

$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $handbackUrl, "", "", $handbackDisplay)
$wsDe.Range("I8").Style = "HyperLink"
